$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value. All values must be written as
# literal text (matching the source inlineStr cells), so numeric-looking
# strings are forced to Text format first and the style is reset afterwards
# to avoid leaving a stray NumberFormat on the cell.
$updates = @(
    ,@("D2", "24.622.78")
    ,@("E2", "  -0.34%  ")
    ,@("D3", "1.694.14")
    ,@("E3", "  -0.08%  ")
    ,@("D4", "1.005")
    ,@("E4", "  +0.34%  ")
    ,@("D5", "314.63")
    ,@("E5", "  -0.51%  ")
    ,@("E6", "  +0.26%  ")
    ,@("D7", "0.3907")
    ,@("E7", "  -1.16%  ")
    ,@("D8", "0.4045")
    ,@("E8", "  -0.47%  ")
    ,@("D9", "1.494")
    ,@("E9", "  +0.45%  ")
    ,@("D10", "1.004")
    ,@("E10", "  +0.30%  ")
    ,@("D11", "53.27")
    ,@("E11", "  +0.33%  ")
    ,@("D12", "0.08740")
    ,@("E12", "  -1.16%  ")
    ,@("D13", "7.636")
    ,@("E13", "  +5.41%  ")
    ,@("E14", "  +3.99%  ")
    ,@("D15", "0.00001356")
    ,@("E15", "  +2.72%  ")
    ,@("D16", "7.967")
    ,@("E16", "  -1.06%  ")
    ,@("D17", "1.690.61")
    ,@("E17", "  -0.40%  ")
    ,@("D18", "98.42")
    ,@("E18", "  -1.65%  ")
    ,@("D19", "0.07104")
    ,@("E19", "  +1.08%  ")
    ,@("D20", "19.81")
    ,@("E20", "  +1.37%  ")
    ,@("D21", "7.311")
    ,@("E21", "  +4.09%  ")
    ,@("D22", "1.006")
    ,@("E22", "  +0.47%  ")
    ,@("D23", "14.27")
    ,@("E23", "  -0.14%  ")
    ,@("D24", "24.610.97")
    ,@("E24", "  -0.33%  ")
    ,@("D25", "3.011")
    ,@("E25", "  -8.08%  ")
    ,@("D26", "2.349")
    ,@("E26", "  -0.70%  ")
    ,@("D27", "22.76")
    ,@("E27", "  -0.29%  ")
    ,@("D28", "162.50")
    ,@("E28", "  -0.75%  ")
    ,@("D29", "8.416")
    ,@("E29", "  +11.08%  ")
    ,@("D30", "137.23")
    ,@("E30", "  +0.82%  ")
    ,@("D31", "5.234")
    ,@("E31", "  +0.97%  ")
    ,@("D32", "1.876.72")
    ,@("E32", "  -0.39%  ")
    ,@("D33", "0.08883")
    ,@("E33", "  +3.25%  ")
    ,@("D34", "7.545")
    ,@("E34", "  +5.63%  ")
    ,@("D35", "1.046")
    ,@("E35", "  -2.10%  ")
    ,@("D36", "1.987")
    ,@("E36", "  +4.59%  ")
    ,@("D37", "0.02929")
    ,@("E37", "  +7.52%  ")
    ,@("D38", "0.2733")
    ,@("E38", "  -0.41%  ")
    ,@("D39", "10.74")
    ,@("E39", "  -5.67%  ")
    ,@("D40", "14.30")
    ,@("E40", "  -0.87%  ")
    ,@("D41", "0.09110")
    ,@("E41", "  -1.42%  ")
    ,@("D42", "0.7879")
    ,@("E42", "  +3.01%  ")
    ,@("D43", "1.465")
    ,@("E43", "  -0.34%  ")
    ,@("D44", "16.88")
    ,@("E44", "  +5.23%  ")
    ,@("D45", "0.7213")
    ,@("E45", "  +0.48%  ")
    ,@("D46", "2.570")
    ,@("E46", "  -0.11%  ")
    ,@("D47", "4.206")
    ,@("E47", "  -0.16%  ")
    ,@("D48", "1.003")
    ,@("E48", "  +0.26%  ")
    ,@("E49", "  -0.02%  ")
    ,@("D50", "138.76")
    ,@("E50", "  -0.50%  ")
    ,@("D51", "91.34")
    ,@("E51", "  +1.41%  ")
)

foreach ($u in $updates) {
    $addr = $u[0]
    $text = $u[1]
    $cell = $ws.Range($addr)
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        # Looks like a plain number (e.g. "1.005") - Excel would otherwise
        # auto-convert the assigned text into a numeric cell. Force the cell
        # to Text first, then drop back to the default "Normal" style so we
        # don't leave a stray number-format behind (source cells carry no
        # explicit style).
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}
